# features.xlsx: for gpio, mention if pull resistor direction has its own register
#
# Add a new "Separate pull resistor reg?" column (I) to the GPIO sheet,
# noting for each MCU whether the pull-resistor direction is controlled
# through its own register or not. Only the first couple of rows are
# filled in (matching the upstream edit), leaving the remaining rows for
# later.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPIO")

# Make GPIO the active/selected sheet (it becomes the front tab on save).
$ws.Select()

$ws.Range("I1").Value = "Separate pull resistor reg?"
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "no, PxOUT"

# Let column widths adjust to fit the new header/content, like Excel does
# automatically as you type into a bestFit column.
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(9).EntireColumn.AutoFit() | Out-Null

# Leave the cursor where the last edit happened.
$ws.Range("I3").Select()
